# Configure the product domain "Assembly" sheet:
#  - insert a new leading column so the existing headers (env..kwargs) shift
#    from B:H to C:I
#  - add a new "Unnamed: 0" header in B1 (the pandas index column header)
#  - add a new "dimensions" header in J1
#  - populate the pandas-style integer index in column A (rows 2-6)
#  - populate the "components" (column D) and "upstream_processes" (column E)
#    data for five assembly rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header columns (env, name, components, upstream_processes,
# downstream_processes, skills, kwargs) one column to the right, B:H -> C:I
$ws.Columns("B").Insert()

# New header cells, matching the bold / bordered / centered style already
# used by the other header cells in row 1
$headerCells = @("B1", "J1")
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("J1").Value = "dimensions"
foreach ($addr in $headerCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Row data: pandas-style integer index (column A) plus component info
$rows = @(
    @{ Row = 2; Index = 0; Name = "casing_2000"; Components = "side_plate_2000: 2, " },
    @{ Row = 3; Index = 1; Name = "casing_3000"; Components = $null },
    @{ Row = 4; Index = 2; Name = "controls";    Components = $null },
    @{ Row = 5; Index = 3; Name = "core_2000";   Components = $null },
    @{ Row = 6; Index = 4; Name = "core_3000";   Components = $null }
)

foreach ($r in $rows) {
    $idxCell = $ws.Cells.Item($r.Row, 1)
    $idxCell.Value = $r.Index
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    $ws.Cells.Item($r.Row, 4).Value = $r.Name

    if ($r.Components) {
        $ws.Cells.Item($r.Row, 5).Value = $r.Components
    }
}

Write-Output "Assembly product domain sheet configured"
